$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Backlog")

# Mark rows (ID 1, 3, 5 -> sheet rows 3, 5, 6) as "Done" and record actual time spent.
$ws.Range("I3").Value = "Done"
$ws.Range("K3").Value = 2

$ws.Range("I5").Value = "Done"
$ws.Range("K5").Value = 1

$ws.Range("I6").Value = "Done"
$ws.Range("K6").Value = 0.5

# Update the view: scroll so column F is the left-most visible column,
# and select I7.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("I7").Select()
